$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("75:75").Insert()

$ws.Range("A75").Value = 7
$ws.Range("B75").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C75").Value = "Ñuble"
$ws.Range("D75").Value = 45161
$ws.Range("E75").Value = 16
$ws.Range("F75").Value = 100112001
$ws.Range("G75").Value = "Berenjena"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 80
$ws.Range("K75").Value = 10000
$ws.Range("L75").Value = 10000
$ws.Range("M75").Value = 10000
$ws.Range("N75").Value = "$/caja 60 unidades"
$ws.Range("O75").Value = "Región de Arica y Parinacota"
$ws.Range("P75").Value = 167
$ws.Range("Q75").Value = 60
$ws.Range("R75").Value = "Hortaliza"
